$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, pushing the former "register_13" summary row
# (name/offset/size/comment) down to row 62.
$ws.Rows(61).Insert()

# --- Row 61: new "reserved" register entry at offset 0x60. rggen
#     renumbers the trailing registers, so this new entry keeps the
#     "register_13" name and the old register_13 becomes register_14. ---
$ws.Range("B57").Copy()
$ws.Range("B61").PasteSpecial(-4122)

$ws.Range("C57").Copy()
$ws.Range("C61:D61").PasteSpecial(-4122)

$ws.Range("E57").Copy()
$ws.Range("E61").PasteSpecial(-4122)

$ws.Range("E11").Copy()
$ws.Range("F61:H61").PasteSpecial(-4122)

$ws.Range("C57").Copy()
$ws.Range("I61").PasteSpecial(-4122)

$ws.Range("J62").Copy()
$ws.Range("J61").PasteSpecial(-4122)
$ws.Range("J61").Borders.Item(9).LineStyle = -4142

$ws.Range("K62").Copy()
$ws.Range("K61").PasteSpecial(-4122)
$ws.Range("K61").Borders.Item(9).LineStyle = -4142

$ws.Range("B61").Value = "register_13"
$ws.Range("C61").Value = "0x60"
$ws.Range("E61").Value = "reserved"

# --- Row 62: the register that used to be register_13 is renumbered
#     register_14 (its offset/size/comment text is unchanged). ---
$ws.Range("B62").Value = "register_14"

# Match the author's final selection/scroll position.
$ws.Range("K62").Select()

Write-Output "ok"
